$wb = $excel.ActiveWorkbook

# --- "Employee" sheet: update dragged/filled values (append "i"), then
# remove the unused template column (E) and blank template rows (5-10)
# that got dragged away to the "Student" sheet below.
$emp = $wb.Worksheets.Item("Employee")

$emp.Range("A2").Value = "Johni"
$emp.Range("C2").Value = "JoshNeo5i"
$emp.Range("D2").Value = "Secret@123i"

$emp.Range("A3").Value = "Maryi"
$emp.Range("B3").Value = "Anni"
$emp.Range("C3").Value = "MarNeo5i"
$emp.Range("D3").Value = "Secret@123i"

$emp.Range("A4").Value = "Davidi"
$emp.Range("B4").Value = "Browni"
$emp.Range("C4").Value = "DavNeo5i"
$emp.Range("D4").Value = "Secret@123i"

$emp.Range("A5:E10").EntireRow.Delete()
$emp.Range("E1:E4").Clear()

# --- "Student" sheet: the blank template grid (rows 3-10, same shape as
# the rows removed from "Employee") was dropped here, then boxed with a
# red rectangle border around the whole block.
$stu = $wb.Worksheets.Item("Student")

$stu.Range("A3:E10").RowHeight = 15.35
$stu.Range("A3:E10").Interior.ColorIndex = 9
$stu.Range("A3:E10").Borders.LineStyle = 0
$stu.Range("A3:E10").BorderAround(1, 2, 10)
